$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# BTEC logo pictures live in the headers (first-page header1.xml, default header2.xml):
# rename image1.jpg -> image2.jpg
$h1 = $sec.Headers.Item(1)
if ($h1.Exists -and $h1.Range.InlineShapes.Count -ge 1) {
    $shape1 = $h1.Range.InlineShapes.Item(1).ConvertToShape()
    $shape1.Name = "image2.jpg"
    $shape1.ConvertToInlineShape()
}

$h2 = $sec.Headers.Item(2)
if ($h2.Exists -and $h2.Range.InlineShapes.Count -ge 1) {
    $shape2 = $h2.Range.InlineShapes.Item(1).ConvertToShape()
    $shape2.Name = "image2.jpg"
    $shape2.ConvertToInlineShape()
}

# Pearson logo pictures live in the footers (first-page footer1.xml, default footer2.xml):
# rename image2.png -> image1.png
$f1 = $sec.Footers.Item(1)
if ($f1.Exists -and $f1.Range.InlineShapes.Count -ge 1) {
    $shape3 = $f1.Range.InlineShapes.Item(1).ConvertToShape()
    $shape3.Name = "image1.png"
    $shape3.ConvertToInlineShape()
}

$f2 = $sec.Footers.Item(2)
if ($f2.Exists -and $f2.Range.InlineShapes.Count -ge 1) {
    $shape4 = $f2.Range.InlineShapes.Item(1).ConvertToShape()
    $shape4.Name = "image1.png"
    $shape4.ConvertToInlineShape()
}

Write-Output "done"
